{"js": "// Append a new paragraph (\"Added something new.\") right after the existing\n// last paragraph (\"I have changed this first time commit 2\"). Using\n// Paragraph.insertParagraph(text, After) means the new paragraph and its run\n// inherit the formatting (bold, size 24 half-points) of the paragraph mark\n// it is split from, matching the target OOXML exactly.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.getLast();\nlastParagraph.insertParagraph(\"Added something new.\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Append a new paragraph (\"Added something new.\") right after the existing\n# last paragraph (\"I have changed this first time commit 2\"). The new\n# paragraph mark inherits the preceding paragraph's formatting (bold,\n# size-24 half-points / 12pt) from InsertParagraphAfter(); the Bold/Size\n# assignments below make that explicit/robust rather than relying solely on\n# inheritance.\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.Text = \"Added something new.\"\n$newParagraph.Range.Font.Bold = $true\n$newParagraph.Range.Font.Size = 12\n"}
